$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column O mirrors the formatting of column N for the header/value rows,
# extending the 2010-2020 series with a 2021 data point.
$ws.Cells.Item(4, 14).Copy()
$ws.Cells.Item(4, 15).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(5, 14).Copy()
$ws.Cells.Item(5, 15).PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Cells.Item(4, 15).Value = 2021
$ws.Cells.Item(5, 15).Value = 1.5020015556876996

# Update the active selection to reflect where the author left off editing
$ws.Range("Q5").Select()
